# ------------------------------------------------------------------
# Applies the tracked edit to QUALIFIED_LIST_OF_STUDENTS.docx:
#   1. Moves the "_GoBack" last-edit bookmark from the "Semester" line
#      to right after the "EE-503" course-code text.
#   2. Recomputes the last table row's monthly attendance numbers
#      (merging what used to be two rows, #11 "21ME19" and
#      #12 "21ME20", into a single corrected row #10 "21ME20") and
#      removes the now-redundant trailing row.
#   3. Removes the stray empty paragraph at the very end of the body.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1a. Remove the existing "_GoBack" bookmark (next to "Spring (2024)") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 1b. Re-create it collapsed, right after the "EE-503" text -------------
# NOTE: a bookmark collapsed exactly on a paragraph's final (mark) position
# cannot be created directly, so we temporarily insert a placeholder
# character after "EE-503", anchor the bookmark there (now a safe,
# non-paragraph-final offset), then delete the placeholder again; the
# collapsed bookmark stays put.
$rng = $d.Content
$found = $rng.Find.Execute("EE-503")
if ($found) {
    $pos = $rng.End
    $placeholder = $d.Range($pos, $pos)
    $placeholder.InsertAfter("X")

    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $d.Range($pos, $pos + 1).Delete()
}

# --- 2. Fix up the table's final rows --------------------------------------
$t = $d.Tables.Item(1)
$rowCount = $t.Rows.Count
$lastRow = $t.Rows.Item($rowCount)
$secondLastRow = $t.Rows.Item($rowCount - 1)

# Merge the two rows' data into the second-to-last row, then drop the last.
$newValues = @("10", "21ME20", "3", "3", "6", "6", "6", "6", "18", "77", "Eligible")
for ($c = 1; $c -le $newValues.Length; $c++) {
    $secondLastRow.Cells.Item($c).Range.Text = $newValues[$c - 1]
}
$lastRow.Delete()

# --- 3. Drop the stray empty paragraph right before the section break ------
$paras = $d.Paragraphs
$count = $paras.Count
$lastPara = $paras.Item($count)
if ($lastPara.Range.Text -eq "\r" -or $lastPara.Range.Text.Trim() -eq "") {
    $prevPara = $paras.Item($count - 1)
    $d.Range($prevPara.Range.End - 1, $lastPara.Range.End).Delete()
}
